# Applies the "Initial commit of login functionality" edit:
#  1. Heading "Authenticating User - Feature Set - 1 Week" gains a tab
#     stop plus a trailing " - (20th-27th February)" + tab.
#  2. Heading "Managing Session History - FS- 2 Weeks ..." gains a
#     trailing " - (27th-13th)" with the _GoBack bookmark now sitting
#     just before the closing paren.
#  3. Heading "Running a Session - FS - 1 week" gains a trailing
#     " - (13th - 22nd March)".
#  4. The two "Will <stop> if ..." bullets become single plain runs.
#  5. The stray _GoBack bookmark that used to sit after "Pull the
#     latest feedback for a current session" is gone (it moved to #2).

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper-free inline approach: build each insertion from a fresh
# zero-length Range so previously-applied character formatting
# (e.g. superscript) never bleeds into the next run.
# ---------------------------------------------------------------

# --- 1) "Authenticating User - Feature Set - 1 Week" heading ----
$p1 = $d.Paragraphs(1)
$p1.TabStops.Add(362.25)          # 7245 twips (362.25 pt * 20)

$r = $p1.Range
$r.End = $r.End - 1               # stop before the paragraph mark
$pos = $r.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" – (20")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("th")
$ins.Font.Superscript = $true
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("-27")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("th")
$ins.Font.Superscript = $true
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" February)")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter([char]9)
$pos = $ins.End

# --- 2) "Managing Session History - FS- 2 Weeks ..." heading ----
$p3 = $d.Paragraphs(3)
$r = $p3.Range
$r.End = $r.End - 1
$pos = $r.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" – (27")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("th")
$ins.Font.Superscript = $true
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("-13")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("th")
$ins.Font.Superscript = $true
$pos = $ins.End

# Insert the closing paren BEFORE adding the bookmark: placing a
# bookmark exactly at a paragraph's end position confuses this
# runtime, so we first give it a normal in-text position to sit at.
$ins = $d.Range($pos, $pos)
$ins.InsertAfter(")")

# _GoBack already exists elsewhere in the doc; re-adding it under the
# same name moves it here (and removes the old one) in one step.
$bm = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bm)

# --- 3) "Running a Session - FS - 1 week" heading ----------------
$p11 = $d.Paragraphs(11)
$r = $p11.Range
$r.End = $r.End - 1
$pos = $r.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" – (13")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("th")
$ins.Font.Superscript = $true
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" – 22")
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("nd")
$ins.Font.Superscript = $true
$pos = $ins.End

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" March)")
$pos = $ins.End

# --- 4) Collapse the two "Will ... stop ... " bullets into a single
#        plain run each. Clearing first (then inserting) forces a
#        real rebuild instead of a same-text no-op. ----------------
$p26 = $d.Paragraphs(26)
$r = $p26.Range
$r.End = $r.End - 1
$r.Text = ""
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertAfter("Will stop if reached lecture end time")

$p27 = $d.Paragraphs(27)
$r = $p27.Range
$r.End = $r.End - 1
$r.Text = ""
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertAfter("Will stop if explicitly told to end early")

# --- 5) The old _GoBack bookmark after "Pull the latest feedback
#        for a current session" is already gone - Bookmarks.Add in
#        step 2 relocated it. Nothing further to do here.
